$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.203.92'
$ws.Range("E2").Value = '  +5.51%  '
$ws.Range("D3").Value = '4.054.72'
$ws.Range("E3").Value = '  +5.71%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'529.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.17%  '
$ws.Range("D6").Value = "'151.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.04%  '
$ws.Range("D7").Value = "'0.710"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +17.80%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = "'0.771"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.80%  '
$ws.Range("E10").Value = '  +7.51%  '
$ws.Range("D11").Value = "'0.0000335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.89%  '
$ws.Range("D12").Value = "'50.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +23.10%  '
$ws.Range("D13").Value = "'11.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.54%  '
$ws.Range("D14").Value = '4.702.55'
$ws.Range("E14").Value = '  +5.88%  '
$ws.Range("D15").Value = '4.052.38'
$ws.Range("E15").Value = '  +5.68%  '
$ws.Range("D16").Value = "'14.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.74%  '
$ws.Range("D17").Value = "'21.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.89%  '
$ws.Range("E18").Value = '  +3.72%  '
$ws.Range("D19").Value = "'0.133"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = '72.174.90'
$ws.Range("E20").Value = '  +5.52%  '
$ws.Range("D21").Value = "'439.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.88%  '
$ws.Range("D22").Value = "'100.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +17.16%  '
$ws.Range("D23").Value = "'3.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.79%  '
$ws.Range("E24").Value = '  +8.80%  '
$ws.Range("D25").Value = "'4.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.83%  '
$ws.Range("D26").Value = "'11.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").Value = "'11.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.78%  '
$ws.Range("D28").Value = "'37.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.89%  '
$ws.Range("E29").Value = '  +3.34%  '
$ws.Range("D30").Value = "'3.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +21.48%  '
$ws.Range("D31").Value = "'13.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.83%  '
$ws.Range("E32").Value = '  +8.01%  '
$ws.Range("D33").Value = "'676.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("D34").Value = "'6.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.63%  '
$ws.Range("D35").Value = "'66.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.46%  '
$ws.Range("D36").Value = "'42.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.73%  '
$ws.Range("D37").Value = "'0.442"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '0.0₃0873'
$ws.Range("E38").Value = '  +6.28%  '
$ws.Range("E39").Value = '  +8.01%  '
$ws.Range("E40").Value = '  +2.80%  '
$ws.Range("D41").Value = "'0.0507"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.06%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").Value = "'3.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.42%  '
$ws.Range("E45").Value = '  +12.53%  '
$ws.Range("D46").Value = "'2.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").Value = "'3.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.03%  '
$ws.Range("D48").Value = "'9.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +14.56%  '
$ws.Range("E49").Value = '  +5.96%  '
$ws.Range("D50").Value = "'3.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.67%  '
$ws.Range("D51").Value = "'0.000276"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.85%  '
